$d = $word.ActiveDocument

# Locate the "12月17日会议" heading paragraph that starts the block of six
# paragraphs (heading + 5 notes) which must collapse into a single empty
# paragraph (keeping only the final paragraph's pPr/rPr formatting).
$findRange = $d.Content
$found = $findRange.Find.Execute("12月17日会议", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate the '12月17日会议' heading paragraph."
}

$targetStart = $findRange.Start

# NOTE: Paragraph.Index is unreliable in this runtime (observed to report a
# value that does not match the actual 1-based position usable with
# $d.Paragraphs.Item(n)). Resolve the real 1-based paragraph number by
# matching the Range.Start offset instead.
$count = $d.Paragraphs.Count
$startIndex = -1
for ($i = 1; $i -le $count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Start -eq $targetStart) {
        $startIndex = $i
        break
    }
}
if ($startIndex -eq -1) {
    throw "Could not resolve the paragraph index for the located heading."
}

# The block runs from the heading paragraph through the end of the document
# (the "4、关于接下来的分工" paragraph is the last paragraph in the body).
$totalParasToMerge = $d.Paragraphs.Count - $startIndex

# Merge the following paragraphs into the heading paragraph by repeatedly
# deleting the paragraph mark right after it. Doing this one paragraph mark
# at a time (rather than deleting one big cross-paragraph range) is what
# this runtime requires to actually collapse the paragraphs. Because the
# merge always happens at the same paragraph mark position, the final
# merged paragraph inherits the identity/formatting (pPr) of the very last
# paragraph in the block, exactly as the target document expects.
for ($i = 0; $i -lt $totalParasToMerge; $i++) {
    $p = $d.Paragraphs.Item($startIndex)
    $markStart = $p.Range.End - 1
    $markRange = $d.Range($markStart, $markStart + 1)
    $markRange.Delete()
}

# Now clear all the run text in the merged paragraph, leaving just the bare
# paragraph mark with the formatting (pPr/rPr rFonts hint) it already has.
$merged = $d.Paragraphs.Item($startIndex)
$textRange = $d.Range($merged.Range.Start, $merged.Range.End - 1)
$textRange.Text = ""
